$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 2.06
$ws.Range("F3").Value = 1.13
$ws.Range("G3").Value = 38
$ws.Range("H3").Value = 1.14
$ws.Range("J3").Value = 7.4
$ws.Range("U3").Value = 1.44
$ws.Range("AE3").Value = 20
$ws.Range("H4").Value = 2.58
$ws.Range("I4").Value = 3.5
$ws.Range("J4").Value = 2.78
$ws.Range("K4").Value = 4.8
$ws.Range("Q4").Value = 1.98
$ws.Range("V4").Value = 1.42
$ws.Range("F5").Value = 7
$ws.Range("G5").Value = 8.6
$ws.Range("H5").Value = 1.46
$ws.Range("I5").Value = 1.53
$ws.Range("J5").Value = 4.7
$ws.Range("K5").Value = 5.2
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 4.5
$ws.Range("P5").Value = 2.18
$ws.Range("Q5").Value = 1.71
$ws.Range("R5").Value = 1.46
$ws.Range("S5").Value = 2.78
$ws.Range("T5").Value = 1.84
$ws.Range("U5").Value = 1.94
$ws.Range("V5").Value = 2.88
$ws.Range("W5").Value = 1.13
$ws.Range("X5").Value = 20
$ws.Range("Y5").Value = 9.800000000000001
$ws.Range("Z5").Value = 9.4
$ws.Range("AA5").Value = 14
$ws.Range("AB5").Value = 990
$ws.Range("AC5").Value = 11.5
$ws.Range("AD5").Value = 11
$ws.Range("AE5").Value = 15.5
$ws.Range("AF5").Value = 70
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 25
$ws.Range("AI5").Value = 980
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 120
$ws.Range("AL5").Value = 110
$ws.Range("AM5").Value = 140
$ws.Range("AN5").Value = 140
$ws.Range("AO5").Value = 7.4
$ws.Range("G6").Value = 1.44
$ws.Range("U6").Value = 1.86
$ws.Range("W6").Value = 3.25
$ws.Range("Z6").Value = 100
$ws.Range("G7").Value = 1.88
$ws.Range("K7").Value = 3.95
$ws.Range("T7").Value = 1.85
$ws.Range("U7").Value = 1.96
$ws.Range("W7").Value = 2.12
$ws.Range("F8").Value = 1.42
$ws.Range("H8").Value = 9.800000000000001
$ws.Range("I8").Value = 10
$ws.Range("L8").Value = 1.37
$ws.Range("O8").Value = 1.31
$ws.Range("V8").Value = 1.11
$ws.Range("Y8").Value = 27
$ws.Range("AA8").Value = 390
$ws.Range("AC8").Value = 11
$ws.Range("AD8").Value = 36
$ws.Range("AE8").Value = 170
$ws.Range("AF8").Value = 7.6
$ws.Range("AH8").Value = 30
$ws.Range("AI8").Value = 150
$ws.Range("AJ8").Value = 11.5
$ws.Range("AN8").Value = 7.2
$ws.Range("AO8").Value = 270
$ws.Range("F9").Value = 2.28
$ws.Range("G9").Value = 2.46
$ws.Range("H9").Value = 3.3
$ws.Range("I9").Value = 3.65
$ws.Range("J9").Value = 3.3
$ws.Range("K9").Value = 3.6
$ws.Range("N9").Value = 3.25
$ws.Range("T9").Value = 1.8
$ws.Range("U9").Value = 2.02
$ws.Range("V9").Value = 1.38
$ws.Range("W9").Value = 1.68
$ws.Range("AB9").Value = 10.5
$ws.Range("J10").Value = 3.5
$ws.Range("K10").Value = 3.6
$ws.Range("P10").Value = 1.78
$ws.Range("Q10").Value = 1.95
$ws.Range("T10").Value = 1.86
$ws.Range("U10").Value = 2
$ws.Range("V10").Value = 1.84
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 7.2
$ws.Range("H11").Value = 1.61
$ws.Range("K11").Value = 4.2
$ws.Range("P11").Value = 1.8
$ws.Range("T11").Value = 2.26
$ws.Range("F12").Value = 2.82
$ws.Range("H12").Value = 2.88
$ws.Range("P12").Value = 1.7
$ws.Range("Q12").Value = 2.36
$ws.Range("T12").Value = 1.97
$ws.Range("U13").Value = 1.89
$ws.Range("G14").Value = 3.7
$ws.Range("H14").Value = 2.4
$ws.Range("K14").Value = 3.2
$ws.Range("AC14").Value = 7
$ws.Range("G15").Value = 2.18
$ws.Range("I15").Value = 4.3
$ws.Range("T15").Value = 1.68
$ws.Range("W15").Value = 1.85
$ws.Range("W16").Value = 1.75
$ws.Range("AC16").Value = 7.6
$ws.Range("I17").Value = 5.6
$ws.Range("Q17").Value = 2.04
$ws.Range("R17").Value = 1.35
$ws.Range("H18").Value = 2.14
$ws.Range("I18").Value = 2.2
$ws.Range("J18").Value = 3.45
$ws.Range("K18").Value = 3.6
$ws.Range("N18").Value = 3.35
$ws.Range("P18").Value = 1.82
$ws.Range("Q18").Value = 2.08
$ws.Range("T18").Value = 1.85
$ws.Range("U18").Value = 2.02
$ws.Range("V18").Value = 1.83
$ws.Range("Z18").Value = 1000
$ws.Range("AC18").Value = 8
$ws.Range("AH18").Value = 1000
$ws.Range("AK18").Value = 50
$ws.Range("AL18").Value = 60
$ws.Range("H19").Value = 2.08
$ws.Range("I19").Value = 2.1
$ws.Range("J19").Value = 3.7
$ws.Range("K19").Value = 3.75
$ws.Range("V19").Value = 1.91
$ws.Range("AH19").Value = 16.5
$ws.Range("I20").Value = 27
$ws.Range("Q20").Value = 1.43
$ws.Range("S20").Value = 2.02
$ws.Range("U20").Value = 1.68
